{"js": "// Replace each three-digit x one-digit multiplication expression with its\n// updated counterpart, per the commit's regenerated problem set.\nconst replacements = [\n  { old: \"149\u00d73=\", new: \"489\u00d78=\" },\n  { old: \"641\u00d73=\", new: \"564\u00d76=\" },\n  { old: \"282\u00d73=\", new: \"832\u00d75=\" },\n  { old: \"671\u00d73=\", new: \"534\u00d72=\" },\n  { old: \"436\u00d72=\", new: \"978\u00d74=\" },\n  { old: \"863\u00d72=\", new: \"889\u00d79=\" },\n  { old: \"168\u00d72=\", new: \"165\u00d74=\" },\n  { old: \"686\u00d77=\", new: \"254\u00d75=\" },\n  { old: \"945\u00d79=\", new: \"775\u00d73=\" },\n  { old: \"693\u00d78=\", new: \"575\u00d78=\" },\n  { old: \"478\u00d78=\", new: \"689\u00d78=\" },\n  { old: \"991\u00d72=\", new: \"478\u00d77=\" },\n  { old: \"380\u00d72=\", new: \"888\u00d76=\" },\n  { old: \"979\u00d79=\", new: \"418\u00d74=\" },\n  { old: \"297\u00d78=\", new: \"181\u00d74=\" },\n  { old: \"155\u00d77=\", new: \"978\u00d74=\" },\n  { old: \"392\u00d79=\", new: \"639\u00d78=\" },\n  { old: \"494\u00d79=\", new: \"621\u00d78=\" },\n  { old: \"550\u00d77=\", new: \"501\u00d72=\" },\n  { old: \"485\u00d78=\", new: \"131\u00d78=\" },\n  { old: \"578\u00d78=\", new: \"521\u00d77=\" },\n  { old: \"837\u00d72=\", new: \"959\u00d79=\" },\n  { old: \"554\u00d72=\", new: \"635\u00d72=\" },\n  { old: \"190\u00d76=\", new: \"524\u00d76=\" },\n  { old: \"324\u00d77=\", new: \"423\u00d77=\" },\n];\n\nconst body = context.document.body;\n\nfor (const { old, new: replacement } of replacements) {\n  const results = body.search(old, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(replacement, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each three-digit x one-digit multiplication expression with its\n# updated counterpart, per the commit's regenerated problem set.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"149\u00d73=\"; New = \"489\u00d78=\" },\n    @{ Old = \"641\u00d73=\"; New = \"564\u00d76=\" },\n    @{ Old = \"282\u00d73=\"; New = \"832\u00d75=\" },\n    @{ Old = \"671\u00d73=\"; New = \"534\u00d72=\" },\n    @{ Old = \"436\u00d72=\"; New = \"978\u00d74=\" },\n    @{ Old = \"863\u00d72=\"; New = \"889\u00d79=\" },\n    @{ Old = \"168\u00d72=\"; New = \"165\u00d74=\" },\n    @{ Old = \"686\u00d77=\"; New = \"254\u00d75=\" },\n    @{ Old = \"945\u00d79=\"; New = \"775\u00d73=\" },\n    @{ Old = \"693\u00d78=\"; New = \"575\u00d78=\" },\n    @{ Old = \"478\u00d78=\"; New = \"689\u00d78=\" },\n    @{ Old = \"991\u00d72=\"; New = \"478\u00d77=\" },\n    @{ Old = \"380\u00d72=\"; New = \"888\u00d76=\" },\n    @{ Old = \"979\u00d79=\"; New = \"418\u00d74=\" },\n    @{ Old = \"297\u00d78=\"; New = \"181\u00d74=\" },\n    @{ Old = \"155\u00d77=\"; New = \"978\u00d74=\" },\n    @{ Old = \"392\u00d79=\"; New = \"639\u00d78=\" },\n    @{ Old = \"494\u00d79=\"; New = \"621\u00d78=\" },\n    @{ Old = \"550\u00d77=\"; New = \"501\u00d72=\" },\n    @{ Old = \"485\u00d78=\"; New = \"131\u00d78=\" },\n    @{ Old = \"578\u00d78=\"; New = \"521\u00d77=\" },\n    @{ Old = \"837\u00d72=\"; New = \"959\u00d79=\" },\n    @{ Old = \"554\u00d72=\"; New = \"635\u00d72=\" },\n    @{ Old = \"190\u00d76=\"; New = \"524\u00d76=\" },\n    @{ Old = \"324\u00d77=\"; New = \"423\u00d77=\" }\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair.Old\n    $find.Replacement.Text = $pair.New\n    $find.Execute([ref]$pair.Old, $true, $true, $false, $false, $false, $true, 1, $false, [ref]$pair.New, 2)\n}\n"}
